# Fertigstellungsgrad Statements.xlsx - add "Write-Statement" completion column formulas
# and update the active selection/view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2 and C3 get their own (non-shared) formula instances.
$ws.Range("C2").Formula = "=IF(B2>0,1,0)"
$ws.Range("C3").Formula = "=IF(B3>0,1,0)"

# C4:C50 are filled in one shot so Excel records them as a shared formula
# group (master formula on C4, followers referencing si="0").
$ws.Range("C4:C50").Formula = "=IF(B4>0,1,0)"

# Move the selection/active cell to E36 (this also clears the custom
# topLeftCell scroll position that was previously stored on the view).
$ws.Range("E36").Select() | Out-Null
